# bird_species_new_add.xlsx — "looked at the newly downloaded species"
#
# Adds a "Migration Period" column (F = autumn table, S = spring table) right
# after the eBird Code column, and removes the per-species "Unique Fact"
# column from the spring table (its text is cleared, but the now-shifted
# empty column is left in place so the used range still spans to column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at C ("Migration Group" and everything right of it
#    shifts one column over, from C:F to D:G).
$ws.Columns("C:C").Insert()

# 2) New column header + values.
$ws.Range("C1").Value = "Migration Period"
$ws.Range("C2:C7").Value = "F"
# C8 (spring-table header row) intentionally stays blank.
$ws.Range("C9:C15").Value = "S"

# 3) The old "Unique Fact" column (now shifted from F to G) is dropped —
#    clear its contents but keep the styled, now-empty cells in place.
$ws.Range("G1:G15").ClearContents()

# 4) Widen the (now) "Why they are in Europe" / notes column F to fit text.
$ws.Columns("F:F").ColumnWidth = 120.8

# 5) Match the author's final selection.
$ws.Range("C15").Select()
